$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    # Force text interpretation so numeric-looking strings (e.g. "605.98")
    # are not auto-converted into floating point numbers, then strip the
    # temporary number-format tweak so the cell's style stays untouched.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.102.30"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.752.19"
$ws.Range("E3").Value = "  +3.72%  "

# Row 5 - BNB
Set-TextValue "D5" "605.98"
$ws.Range("E5").Value = "  +1.31%  "

# Row 6 - Solana
Set-TextValue "D6" "167.68"
$ws.Range("E6").Value = "  +5.23%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.751.38"
$ws.Range("E9").Value = "  +3.75%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.92%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.43%  "

# Row 14 - Avalanche
Set-TextValue "D14" "28.96"
$ws.Range("E14").Value = "  +2.90%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.252.36"
$ws.Range("E15").Value = "  +3.72%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -0.14%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.025.79"
$ws.Range("E17").Value = "  +1.13%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.756.11"
$ws.Range("E18").Value = "  +4.10%  "

# Row 19 - Chainlink
Set-TextValue "D19" "12.03"
$ws.Range("E19").Value = "  +5.11%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.75"
$ws.Range("E20").Value = "  +5.34%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "369.80"
$ws.Range("E21").Value = "  +1.32%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  +3.17%  "

# Row 23 - NEARProtocol
$ws.Range("E23").Value = "  +3.20%  "

# Row 24 - SuiNetwork
Set-TextValue "D24" "2.14"
$ws.Range("E24").Value = "  +3.06%  "

# Row 25 - Litecoin
Set-TextValue "D25" "74.23"
$ws.Range("E25").Value = "  -1.23%  "

# Row 27 - Aptos
Set-TextValue "D27" "10.00"
$ws.Range("E27").Value = "  +2.40%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +1.46%  "

# Row 30 - Bittensor
Set-TextValue "D30" "602.02"
$ws.Range("E30").Value = "  +7.42%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  -1.88%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +3.91%  "

# Row 33 - Fetch.AI
Set-TextValue "D33" "1.47"

# Row 34 - PancakeSwap
Set-TextValue "D34" "1.98"
$ws.Range("E34").Value = "  +6.04%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +3.52%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +4.16%  "

# Row 37 - FirstDigitalUSD
Set-TextValue "D37" "0.999"
$ws.Range("E37").Value = "  -0.03%  "

# Row 38 & 39 - swap Monero and EthereumClassic, with new values
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D38" "20.24"
$ws.Range("E38").Value = "  +1.77%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D39" "163.36"
$ws.Range("E39").Value = "  +2.34%  "

# Row 40 - PolygonEcosystemToken
Set-TextValue "D40" "0.385"
$ws.Range("E40").Value = "  +3.55%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  +2.52%  "

# Row 42 - RenderToken
Set-TextValue "D42" "5.53"
$ws.Range("E42").Value = "  +3.02%  "

# Row 43 - dogwifhat
Set-TextValue "D43" "2.72"
$ws.Range("E43").Value = "  +3.28%  "

# Row 44 - WhiteBITCoin
$ws.Range("E44").Value = "  +1.30%  "

# Row 45 - BabyDogeCoin
$ws.Range("E45").Value = "  -5.04%  "

# Row 47 - Aave
Set-TextValue "D47" "159.36"
$ws.Range("E47").Value = "  +0.67%  "

# Row 48 - Filecoin
$ws.Range("E48").Value = "  +5.35%  "

# Row 49 - Optimism
$ws.Range("E49").Value = "  +6.75%  "

# Row 50 - ARBITRUM
$ws.Range("E50").Value = "  +7.80%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "22.19"
$ws.Range("E51").Value = "  -0.40%  "
